$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Productivity [MMGGE/yr]" column (old column H) entirely;
# this shifts I/J/K left to H/I/J and updates the dimension/merged range.
$ws.Columns("H").Delete()

# Refresh data values (new spearman run) across the remaining value columns.
$ws.Range("C4").Value = 0.02622262226222622
$ws.Range("E4").Value = 0.05513351335133513
$ws.Range("F4").Value = 0.07425142514251425
$ws.Range("H4").Value = 0.06303030303030302
$ws.Range("I4").Value = 0.06549054905490548
$ws.Range("J4").Value = -0.06365836848766269
$ws.Range("C5").Value = 0.1291689168916892
$ws.Range("E5").Value = 0.0003000300030002999
$ws.Range("F5").Value = 0.007116711671167116
$ws.Range("H5").Value = -0.06853885388538854
$ws.Range("I5").Value = 0.007644764476447644
$ws.Range("J5").Value = 0.09616298104940993
$ws.Range("C6").Value = -0.04182418241824182
$ws.Range("E6").Value = -0.07277527752775277
$ws.Range("F6").Value = -0.08244824482448244
$ws.Range("H6").Value = -0.01870987098709871
$ws.Range("I6").Value = -0.0738073807380738
$ws.Range("J6").Value = 0.006273978597529711
$ws.Range("C7").Value = 0.03228322832283228
$ws.Range("E7").Value = 0.7706330633063305
$ws.Range("F7").Value = 0.7708610861086107
$ws.Range("H7").Value = 0.9997839783978397
$ws.Range("I7").Value = 0.7682568256825681
$ws.Range("J7").Value = -0.1074561425249634
$ws.Range("C8").Value = 0.9654965496549653
$ws.Range("E8").Value = 0.05034503450345034
$ws.Range("F8").Value = 0.04788478847884788
$ws.Range("H8").Value = 0.01663366336633663
$ws.Range("I8").Value = 0.0403000300030003
$ws.Range("J8").Value = 0.1220453941919464
$ws.Range("C9").Value = 0.134029402940294
$ws.Range("E9").Value = 0.1125112511251125
$ws.Range("F9").Value = 0.1123072307230723
$ws.Range("H9").Value = 0.121980198019802
$ws.Range("I9").Value = 0.1078307830783078
$ws.Range("J9").Value = 0.1595031477517434
$ws.Range("C10").Value = 0.00522052205220522
$ws.Range("E10").Value = 0.1051905190519052
$ws.Range("F10").Value = 0.1362976297629763
$ws.Range("H10").Value = -0.02341434143414341
$ws.Range("I10").Value = 0.123048304830483
$ws.Range("J10").Value = -0.02253829249294405
$ws.Range("C11").Value = 0.1054065406540654
$ws.Range("E11").Value = -0.1616921692169217
$ws.Range("F11").Value = -0.1412901290129013
$ws.Range("H11").Value = -0.1034623462346234
$ws.Range("I11").Value = -0.1527512751275127
$ws.Range("J11").Value = -0.2168575186055245
$ws.Range("C12").Value = 0.1404740474047405
$ws.Range("E12").Value = 0.4869126912691268
$ws.Range("F12").Value = 0.4854965496549655
$ws.Range("H12").Value = -0.1544554455445544
$ws.Range("I12").Value = 0.4959975997599759
$ws.Range("J12").Value = 0.03774593630877444
$ws.Range("C13").Value = -0.2772277227722772
$ws.Range("E13").Value = -0.008076807680768076
$ws.Range("F13").Value = 0.005376537653765376
$ws.Range("H13").Value = 0.00534053405340534
$ws.Range("I13").Value = -0.004776477647764776
$ws.Range("J13").Value = -0.1510858099586939
